$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.304.79'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '3.267.19'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '3.268.84'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.81'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.394'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').Value = '3.836.14'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').Value = '65.455.05'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '3.281.96'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000160'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '418.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.501'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.204'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000111'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.30'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.64%  '
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.16'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').Value = '2.829.39'
$ws.Range('E39').Value = '  +2.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.739'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('E45').Value = '  -3.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0630'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '310.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('E51').Value = '  -0.84%  '
